$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.981.63'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.237.09'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'306.25"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.10%  '
$ws.Range('D6').Value = "'94.61"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.80%  '
$ws.Range('D7').Value = "'0.569"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('E9').Value = '  -4.56%  '
$ws.Range('E10').Value = '  -5.48%  '
$ws.Range('D11').Value = "'0.0806"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.38%  '
$ws.Range('E12').Value = '  -4.09%  '
$ws.Range('E13').Value = '  -1.34%  '
$ws.Range('D14').Value = '2.576.69'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').Value = '2.235.51'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').Value = "'13.58"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.83%  '
$ws.Range('D18').Value = '43.855.22'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = '0.0₃0959'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('D20').Value = "'12.09"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.19%  '
$ws.Range('E21').Value = '  -2.66%  '
$ws.Range('D22').Value = "'65.01"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = "'236.31"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('E24').Value = '  -5.35%  '
$ws.Range('E25').Value = '  -5.58%  '
$ws.Range('D27').Value = "'9.95"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.67%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = "'2.17"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = "'37.38"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.46%  '
$ws.Range('D30').Value = "'5.98"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').Value = "'19.89"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = "'152.91"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('D33').Value = "'0.0801"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.81%  '
$ws.Range('D34').Value = "'3.26"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.93%  '
$ws.Range('D35').Value = "'2.58"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.62%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -4.93%  '
$ws.Range('E38').Value = '  -7.70%  '
$ws.Range('D39').Value = "'15.13"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.93%  '
$ws.Range('E40').Value = '  -7.48%  '
$ws.Range('E41').Value = '  -8.70%  '
$ws.Range('E42').Value = '  -4.48%  '
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').Value = '1.726.53'
$ws.Range('E44').Value = '  -2.03%  '
$ws.Range('D45').Value = "'85.44"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.55%  '
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').Value = "'99.99"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.04%  '
$ws.Range('D48').Value = "'4.91"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.77%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = "'69.04"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.15%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = "'8.08"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('D51').Value = "'54.09"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.37%  '
